$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update re-associates odds/result data with the correct match id for
# several fixtures. In the sheet, this shows up as blocks of row content
# (columns B..AB) being swapped or cyclically rotated between rows while
# the row-rank column (A) stays put.

# Rows 17 and 20: swap full data block (B:AB)
$r17 = $ws.Range("B17:AB17").Value()
$r20 = $ws.Range("B20:AB20").Value()
$ws.Range("B17:AB17").Value = $r20
$ws.Range("B20:AB20").Value = $r17

# Rows 27 and 28: swap full data block (B:AB)
$r27 = $ws.Range("B27:AB27").Value()
$r28 = $ws.Range("B28:AB28").Value()
$ws.Range("B27:AB27").Value = $r28
$ws.Range("B28:AB28").Value = $r27

# Rows 29, 30, 31: cyclic rotation -> content of 29 goes to 31,
# content of 30 goes to 29, content of 31 goes to 30
$r29 = $ws.Range("B29:AB29").Value()
$r30 = $ws.Range("B30:AB30").Value()
$r31 = $ws.Range("B31:AB31").Value()
$ws.Range("B29:AB29").Value = $r30
$ws.Range("B30:AB30").Value = $r31
$ws.Range("B31:AB31").Value = $r29

# Rows 109 and 110: swap full data block (B:AB)
$r109 = $ws.Range("B109:AB109").Value()
$r110 = $ws.Range("B110:AB110").Value()
$ws.Range("B109:AB109").Value = $r110
$ws.Range("B110:AB110").Value = $r109
